# Update "想去人数" (F column) figures across the sheets to the freshly
# scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (Exhibitions) --------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1304
$ws.Range("F4").Value = 1100
$ws.Range("F5").Value = 985
$ws.Range("F6").Value = 1763
$ws.Range("F7").Value = 543
$ws.Range("F11").Value = 123
$ws.Range("F12").Value = 274
$ws.Range("F13").Value = 58
$ws.Range("F14").Value = 84
$ws.Range("F15").Value = 659
$ws.Range("F16").Value = 149
$ws.Range("F17").Value = 96
$ws.Range("F20").Value = 324
$ws.Range("F21").Value = 126
$ws.Range("F22").Value = 655
$ws.Range("F23").Value = 23
$ws.Range("F26").Value = 35
$ws.Range("F27").Value = 860
$ws.Range("F28").Value = 306
$ws.Range("F29").Value = 143
$ws.Range("F31").Value = 262
$ws.Range("F32").Value = 11

# -- Sheet "演出" (Performances) --------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 10
$ws.Range("F7").Value = 246
$ws.Range("F10").Value = 618

# -- Sheet "全部类型" (All types) --------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1304
$ws.Range("F5").Value = 1100
$ws.Range("F6").Value = 985
$ws.Range("F7").Value = 1763
$ws.Range("F8").Value = 543
$ws.Range("F13").Value = 123
$ws.Range("F14").Value = 274
$ws.Range("F15").Value = 58
$ws.Range("F16").Value = 84
$ws.Range("F17").Value = 659
$ws.Range("F18").Value = 149
$ws.Range("F19").Value = 96
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 324
$ws.Range("F27").Value = 246
$ws.Range("F28").Value = 246
$ws.Range("F29").Value = 126
$ws.Range("F30").Value = 655
$ws.Range("F31").Value = 23
$ws.Range("F34").Value = 35
$ws.Range("F35").Value = 860
$ws.Range("F36").Value = 306
$ws.Range("F39").Value = 143
$ws.Range("F41").Value = 262
$ws.Range("F42").Value = 618
$ws.Range("F45").Value = 11
